$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 49, shifting existing rows 49:64 down to 50:65.
$ws.Rows.Item(49).Insert()

# The newly inserted row 49 is blank; copy static column values and
# formatting from the row immediately below (which now holds what used
# to be row 49) so the new record matches the sheet's formatting.
$src = $ws.Rows.Item(50)
$dst = $ws.Rows.Item(49)
$src.Copy()
$dst.PasteSpecial(-4104) # xlPasteAll

# Now overwrite the new row's values with the new weekly record.
$ws.Cells.Item(49, 1).Value = 7
$ws.Cells.Item(49, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(49, 3).Value = "Ñuble"
$ws.Cells.Item(49, 4).Value = 44588
$ws.Cells.Item(49, 5).Value = 16
$ws.Cells.Item(49, 6).Value = 100112031
$ws.Cells.Item(49, 7).Value = "Poroto verde"
$ws.Cells.Item(49, 8).Value = "Sin especificar"
$ws.Cells.Item(49, 9).Value = "Primera"
$ws.Cells.Item(49, 10).Value = 120
$ws.Cells.Item(49, 11).Value = 24000
$ws.Cells.Item(49, 12).Value = 25000
$ws.Cells.Item(49, 13).Value = 24500
$ws.Cells.Item(49, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(49, 15).Value = "Región del Maule"
$ws.Cells.Item(49, 16).Value = 980
$ws.Cells.Item(49, 17).Value = 25
$ws.Cells.Item(49, 18).Value = "Hortaliza"
